$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ativação: 01/01/2020 -> 01/01/2022
# Pre-format as Text so Excel stores the literal string instead of
# auto-converting the date-like text into a date serial value.
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "01/01/2022"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "01/01/2022"

# Programa resumido: remove "Ciclo Hidrológico; " and trailing period
$ws.Range("B15").Value = "Geomorfologia Fluvial; Padrões de Drenagem; Escoamentos hidráulicos; medidores; bocais; instrumentos de medição"
$ws.Range("C15").Value = "Geomorfologia Fluvial; Padrões de Drenagem; Escoamentos hidráulicos; medidores; bocais; instrumentos de medição"

# Short syllabus: remove "Hydrological Cycle; "
$ws.Range("B16").Value = "River Geomorphology; Drainage Patterns; Hydraulic flow; meters; nozzles; measuring instruments."
$ws.Range("C16").Value = "River Geomorphology; Drainage Patterns; Hydraulic flow; meters; nozzles; measuring instruments."

# Programa: full rewrite to bullet-style list
$ws.Range("B17").Value = "- As teorias geomorfológicas;- Processos e Formas do relevo;- Processos fluviais, morfologias fluviais e padrões de drenagem;- Precipitação;- Infiltração;- Evapotranspiração;- Escoamento superficial;- Instrumentos de medição (Calhas, vertedores e registros);- Operação de reservatórios;- Vazões máximas e mínimas: distribuição de frequência, hidrograma unitário.- Água subterrânea, aquíferos e poços;"
$ws.Range("C17").Value = "- As teorias geomorfológicas;- Processos e Formas do relevo;- Processos fluviais, morfologias fluviais e padrões de drenagem;- Precipitação;- Infiltração;- Evapotranspiração;- Escoamento superficial;- Instrumentos de medição (Calhas, vertedores e registros);- Operação de reservatórios;- Vazões máximas e mínimas: distribuição de frequência, hidrograma unitário.- Água subterrânea, aquíferos e poços;"

# Syllabus: full rewrite to bullet-style list
$ws.Range("B18").Value = "- Geomorphological theories;- Processes and Forms of relief;- River processes, river morphologies and drainage patterns;- Precipitation;- Infiltration;- Evapotranspiration;- Surface runoff;- Measuring instruments (gutters, spillways and registers);- Reservoir operation;- Maximum and minimum flow rates: frequency distribution, unit hydrograph.- Groundwater, aquifers and wells;"
$ws.Range("C18").Value = "- Geomorphological theories;- Processes and Forms of relief;- River processes, river morphologies and drainage patterns;- Precipitation;- Infiltration;- Evapotranspiration;- Surface runoff;- Measuring instruments (gutters, spillways and registers);- Reservoir operation;- Maximum and minimum flow rates: frequency distribution, unit hydrograph.- Groundwater, aquifers and wells;"
